$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.8520199728707212
$ws.Range("C2").Value = 0.2234733660233132
$ws.Range("D2").Value = 0.0797080734427027
$ws.Range("E2").Value = 0.425497448737346
$ws.Range("G2").Value = 0.3167004636070629
$ws.Range("H2").Value = 0.4419693898604891
$ws.Range("N2").Value = 0.8079598454597203
$ws.Range("O2").Value = 1.446305228668251
$ws.Range("B3").Value = 0.7464385980046018
$ws.Range("C3").Value = 0.1985273790015185
$ws.Range("D3").Value = 0.07213913390636151
$ws.Range("E3").Value = 0.3710987126234215
$ws.Range("G3").Value = 0.3061705257600238
$ws.Range("H3").Value = 0.4415624552119937
$ws.Range("N3").Value = 0.812434400066671
$ws.Range("O3").Value = 1.422795646831162
$ws.Range("B4").Value = 0.6815076707742946
$ws.Range("C4").Value = 0.1831284550079317
$ws.Range("D4").Value = 0.06752681047800024
$ws.Range("E4").Value = 0.3377923118596016
$ws.Range("G4").Value = 0.300034475101782
$ws.Range("H4").Value = 0.4416194215679496
$ws.Range("N4").Value = 0.8155669144548767
$ws.Range("O4").Value = 1.409689370933592
$ws.Range("B5").Value = 0.6550227555186154
$ws.Range("C5").Value = 0.1768329077078761
$ws.Range("D5").Value = 0.0656560306131837
$ws.Range("E5").Value = 0.3242415683853039
$ws.Range("G5").Value = 0.2976162070390131
$ws.Range("H5").Value = 0.4417196117754543
$ws.Range("N5").Value = 0.8169403745251529
$ws.Range("O5").Value = 1.404680914476586
$ws.Range("B6").Value = 0.6506234731130007
$ws.Range("C6").Value = 0.175786314713605
$ws.Range("D6").Value = 0.06534591841133874
$ws.Range("E6").Value = 0.321992743883655
$ws.Range("G6").Value = 0.2972196039344226
$ws.Range("H6").Value = 0.441740892262473
$ws.Range("N6").Value = 0.8171742940355671
$ws.Range("O6").Value = 1.403869295376779
$ws.Range("B7").Value = 0.6811505859801059
$ws.Range("C7").Value = 0.183043633153261
$ws.Range("D7").Value = 0.06750154497855476
$ws.Range("E7").Value = 0.3376094756890495
$ws.Range("G7").Value = 0.3000015293748675
$ws.Range("H7").Value = 0.4416204613336987
$ws.Range("N7").Value = 0.8155850448071078
$ws.Range("O7").Value = 1.409620481231372
$ws.Range("B8").Value = 0.8156375692544202
$ws.Range("C8").Value = 0.2148891515977311
$ws.Range("D8").Value = 0.07709101042431143
$ws.Range("E8").Value = 0.4067199514648792
$ws.Range("G8").Value = 0.3130010133544943
$ws.Range("H8").Value = 0.4417652881954695
$ws.Range("N8").Value = 0.8094227936367133
$ws.Range("O8").Value = 1.437922343378176
$ws.Range("B9").Value = 1.078513806947512
$ws.Range("C9").Value = 0.2766793140977484
$ws.Range("D9").Value = 0.09617631966256113
$ws.Range("E9").Value = 0.5430890381526723
$ws.Range("G9").Value = 0.3411352589640302
$ws.Range("H9").Value = 0.4444925873505241
$ws.Range("N9").Value = 0.8003910041157383
$ws.Range("O9").Value = 1.504044797490877
$ws.Range("B10").Value = 1.27110733669673
$ws.Range("C10").Value = 0.3216681755163222
$ws.Range("D10").Value = 0.1103740216886138
$ws.Range("E10").Value = 0.6439312962175165
$ws.Range("G10").Value = 0.3634587668815783
$ws.Range("H10").Value = 0.447999044992315
$ws.Range("N10").Value = 0.7956126648292496
$ws.Range("O10").Value = 1.559218994920172
$ws.Range("B11").Value = 1.358602154837058
$ws.Range("C11").Value = 0.3420449685987137
$ws.Range("D11").Value = 0.1168721984325742
$ws.Range("E11").Value = 0.689979021072773
$ws.Range("G11").Value = 0.3739826124665626
$ws.Range("H11").Value = 0.4499234614084884
$ws.Range("N11").Value = 0.7938416253101366
$ws.Range("O11").Value = 1.585777000038973
$ws.Range("B12").Value = 1.391716640454945
$ws.Range("C12").Value = 0.34974815995227
$ws.Range("D12").Value = 0.1193386365585383
$ws.Range("E12").Value = 0.7074434256330022
$ws.Range("G12").Value = 0.3780214508763322
$ws.Range("H12").Value = 0.4506997646181361
$ws.Range("N12").Value = 0.7932288370498668
$ws.Range("O12").Value = 1.596045590715363
$ws.Range("B13").Value = 1.384585656927527
$ws.Range("C13").Value = 0.3480897250406088
$ws.Range("D13").Value = 0.1188071901275976
$ws.Range("E13").Value = 0.7036809133413868
$ws.Range("G13").Value = 0.3771492165402748
$ws.Range("H13").Value = 0.4505304547858202
$ws.Range("N13").Value = 0.793358238841904
$ws.Range("O13").Value = 1.593824623218296
$ws.Range("B14").Value = 1.361326868493848
$ws.Range("C14").Value = 0.3426789780783679
$ws.Range("D14").Value = 0.1170749989829716
$ws.Range("E14").Value = 0.6914152723078075
$ws.Range("G14").Value = 0.3743138104001105
$ws.Range("H14").Value = 0.449986373766464
$ws.Range("N14").Value = 0.7937900513061606
$ws.Range("O14").Value = 1.586617550756102
$ws.Range("B15").Value = 1.347077825393967
$ws.Range("C15").Value = 0.3393630305585873
$ws.Range("D15").Value = 0.1160147276281833
$ws.Range("E15").Value = 0.6839058088808798
$ws.Range("G15").Value = 0.3725840533953004
$ws.Range("H15").Value = 0.4496593096293822
$ws.Range("N15").Value = 0.7940620838497381
$ws.Range("O15").Value = 1.582230628967352
$ws.Range("B16").Value = 1.265386898229963
$ws.Range("C16").Value = 0.3203346920404329
$ws.Range("D16").Value = 0.1099501493337414
$ws.Range("E16").Value = 0.6409256508423624
$ws.Range("G16").Value = 0.3627784855973033
$ws.Range("H16").Value = 0.4478799246857079
$ws.Range("N16").Value = 0.7957365056688417
$ws.Range("O16").Value = 1.557512883320413
$ws.Range("B17").Value = 1.215241431867071
$ws.Range("C17").Value = 0.3086384626095366
$ws.Range("D17").Value = 0.1062398829098754
$ws.Range("E17").Value = 0.6146047871471154
$ws.Range("G17").Value = 0.3568580055269592
$ws.Range("H17").Value = 0.4468728258510453
$ws.Range("N17").Value = 0.7968668109276962
$ws.Range("O17").Value = 1.542724342118788
$ws.Range("B18").Value = 1.186388142893861
$ws.Range("C18").Value = 0.3019027528762876
$ws.Range("D18").Value = 0.1041095586376173
$ws.Range("E18").Value = 0.5994819806032723
$ws.Range("G18").Value = 0.3534873728815455
$ws.Range("H18").Value = 0.4463245526658142
$ws.Range("N18").Value = 0.7975548336651599
$ws.Range("O18").Value = 1.534355579889279
$ws.Range("B19").Value = 1.176617068809662
$ws.Range("C19").Value = 0.299620732080399
$ws.Range("D19").Value = 0.1033889059499842
$ws.Range("E19").Value = 0.5943643868934174
$ws.Range("G19").Value = 0.3523520693379254
$ws.Range("H19").Value = 0.4461442317336264
$ws.Range("N19").Value = 0.7977942970585801
$ws.Range("O19").Value = 1.53154557146658
$ws.Range("B20").Value = 1.220580642831237
$ws.Range("C20").Value = 0.3098844122096693
$ws.Range("D20").Value = 0.1066344615928045
$ws.Range("E20").Value = 0.6174049902606811
$ws.Range("G20").Value = 0.3574846576700708
$ws.Range("H20").Value = 0.4469768250250752
$ws.Range("N20").Value = 0.7967425656561034
$ws.Range("O20").Value = 1.544284392898078
$ws.Range("B21").Value = 1.36815903295701
$ws.Range("C21").Value = 0.3442686022216606
$ws.Range("D21").Value = 0.1175836300477897
$ws.Range("E21").Value = 0.695017235251882
$ws.Range("G21").Value = 0.3751451756250077
$ws.Range("H21").Value = 0.4501448908422105
$ws.Range("N21").Value = 0.793661647295508
$ws.Range("O21").Value = 1.588728683719438
$ws.Range("B22").Value = 1.464505191610669
$ws.Range("C22").Value = 0.3666644186203882
$ws.Range("D22").Value = 0.1247728966167472
$ws.Range("E22").Value = 0.74590045162509
$ws.Range("G22").Value = 0.3870005778549626
$ws.Range("H22").Value = 0.4524927709863391
$ws.Range("N22").Value = 0.7919853648367194
$ws.Range("O22").Value = 1.619010122788836
$ws.Range("B23").Value = 1.413093383856221
$ws.Range("C23").Value = 0.3547184174044844
$ws.Range("D23").Value = 0.1209327883019142
$ws.Range("E23").Value = 0.7187278618635133
$ws.Range("G23").Value = 0.3806442477987275
$ws.Range("H23").Value = 0.4512142130389662
$ws.Range("N23").Value = 0.7928491773190132
$ws.Range("O23").Value = 1.602734776216437
$ws.Range("B24").Value = 1.21816686020901
$ws.Range("C24").Value = 0.3093211538622711
$ws.Range("D24").Value = 0.1064560639863004
$ws.Range("E24").Value = 0.6161389893273395
$ws.Range("G24").Value = 0.3572012450858324
$ws.Range("H24").Value = 0.4469297113133166
$ws.Range("N24").Value = 0.7967986179577338
$ws.Range("O24").Value = 1.543578678634987
$ws.Range("B25").Value = 1.007492380989447
$ws.Range("C25").Value = 0.2600346739734505
$ws.Range("D25").Value = 0.09098272485020686
$ws.Range("E25").Value = 0.5060931205242838
$ws.Range("G25").Value = 0.333236942504854
$ws.Range("H25").Value = 0.443491849494194
$ws.Range("N25").Value = 0.8025079528096057
$ws.Range("O25").Value = 1.485007509827227
